$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) cells keep exact text formatting (no numeric auto-conversion)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.596.54"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.798.75"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.63"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("E6").Value = "  +1.64%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.63"
$ws.Range("E8").Value = "  +1.95%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.297"
$ws.Range("E9").Value = "  +2.11%  "

$ws.Range("E10").Value = "  +0.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.059.98"
$ws.Range("E12").Value = "  +0.88%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.17"
$ws.Range("E13").Value = "  +0.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.796.77"
$ws.Range("E14").Value = "  +0.65%  "

$ws.Range("E15").Value = "  +2.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.587.83"
$ws.Range("E16").Value = "  +1.17%  "

$ws.Range("E17").Value = "  +2.90%  "

$ws.Range("E18").Value = "  +1.44%  "

$ws.Range("E19").Value = "  -0.66%  "

$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.31"
$ws.Range("E21").Value = "  +2.90%  "

$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "169.24"
$ws.Range("E24").Value = "  +4.05%  "

$ws.Range("E25").Value = "  +1.45%  "

$ws.Range("E26").Value = "  +1.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.61"
$ws.Range("E27").Value = "  +1.79%  "

$ws.Range("E28").Value = "  +2.22%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.10"
$ws.Range("E30").Value = "  +10.48%  "

$ws.Range("E31").Value = "  +2.42%  "

$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("E33").Value = "  +0.80%  "

$ws.Range("E34").Value = "  +2.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.432.97"
$ws.Range("E35").Value = "  -0.95%  "

$ws.Range("E36").Value = "  +7.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.676"
$ws.Range("E37").Value = "  +3.32%  "

$ws.Range("E38").Value = "  +2.71%  "

$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.13"
$ws.Range("E40").Value = "  +5.84%  "

$ws.Range("E41").Value = "  +1.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.939"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("E43").Value = "  +3.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.96"
$ws.Range("E44").Value = "  +2.24%  "

$ws.Range("E45").Value = "  +3.39%  "

$ws.Range("E46").Value = "  +0.53%  "

$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.959.86"
$ws.Range("E48").Value = "  +0.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.08"
$ws.Range("E49").Value = "  +1.18%  "

$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0129"
$ws.Range("E51").Value = "  -4.35%  "
